$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("M9").ClearContents()
$ws.Range("N9").ClearContents()
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("M18").ClearContents()
$ws.Range("H41").Value = 281.66666
$ws.Range("I41").Value = 96.333336
$ws.Range("J41").Value = 467
$ws.Range("K41").Value = 96.333336
$ws.Range("L41").Value = 467
$ws.Range("M41").Value = 343.666664
$ws.Range("N41").Value = -1347
$ws.Range("H137").Value = 3866
$ws.Range("I137").Value = 3866
$ws.Range("K137").Value = 11598
$ws.Range("M137").Value = -9048
$ws.Range("H138").Value = 2681.818
$ws.Range("J138").Value = 2650
$ws.Range("L138").Value = 7950
$ws.Range("N138").Value = -18230
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 150
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 150
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 150
$ws.Range("M5").ClearContents()
$ws.Range("N5").Value = -374
$ws.Range("H32").Value = 14208.5
$ws.Range("I32").Value = 12224.615
$ws.Range("K32").Value = 12224.615
$ws.Range("M32").Value = -11937.615
$ws.Range("H122").Value = 988.5
$ws.Range("I122").Value = 988.5
$ws.Range("K122").Value = 2965.5
$ws.Range("M122").Value = -515.5
$ws.Range("H132").Value = 3762.375
$ws.Range("I132").Value = 3762.375
$ws.Range("K132").Value = 11287.125
$ws.Range("M132").Value = -8757.125
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 150
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 150
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 150
$ws.Range("M4").ClearContents()
$ws.Range("N4").Value = -380
$ws.Range("H37").Value = 1250
$ws.Range("I37").Value = 1250
$ws.Range("K37").Value = 1250
$ws.Range("M37").Value = -1113
$ws.Range("H134").Value = 1058.7142
$ws.Range("I134").Value = 1058.7142
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 3176.1426
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -641.1425999999997
$ws.Range("N134").ClearContents()
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 30.846153
$ws.Range("I7").Value = 9.142858
$ws.Range("K7").Value = 9.142858
$ws.Range("M7").Value = 103.857142
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("H31").Value = 2879.25
$ws.Range("I31").Value = 1965.3334
$ws.Range("J31").Value = 5621
$ws.Range("K31").Value = 1965.3334
$ws.Range("L31").Value = 5621
$ws.Range("M31").Value = -1670.3334
$ws.Range("N31").Value = -6211
$ws.Range("H34").Value = 2879.25
$ws.Range("I34").Value = 1965.3334
$ws.Range("J34").Value = 5621
$ws.Range("K34").Value = 1965.3334
$ws.Range("L34").Value = 5621
$ws.Range("M34").Value = -1763.3334
$ws.Range("N34").Value = -6025
$ws.Range("H95").Value = 9279.799999999999
$ws.Range("J95").Value = 9279.799999999999
$ws.Range("L95").Value = 9279.799999999999
$ws.Range("N95").Value = -14771.8
$ws.Range("H134").Value = 549.75
$ws.Range("J134").Value = 900
$ws.Range("L134").Value = 2700
$ws.Range("N134").Value = -7770
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 624
$ws.Range("I68").Value = 471.2857
$ws.Range("J68").Value = 837.8
$ws.Range("K68").Value = 1413.8571
$ws.Range("L68").Value = 2513.4
$ws.Range("M68").Value = -602.8571000000002
$ws.Range("N68").Value = -4135.4
$ws.Range("H71").Value = 624
$ws.Range("I71").Value = 471.2857
$ws.Range("J71").Value = 837.8
$ws.Range("K71").Value = 4241.571300000001
$ws.Range("L71").Value = 7540.2
$ws.Range("M71").Value = -185.5713000000005
$ws.Range("N71").Value = -15652.2
$ws.Range("H75").Value = 200
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()
$ws.Range("H78").Value = 200
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()
$ws.Range("H88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("N88").ClearContents()
$ws.Range("H91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("N91").ClearContents()
$ws.Range("H107").Value = 971.3158
$ws.Range("J107").Value = 619.8889
$ws.Range("L107").Value = 1859.6667
$ws.Range("N107").Value = -5699.6667
$ws.Range("H108").Value = 571
$ws.Range("I108").Value = 571
$ws.Range("K108").Value = 1713
$ws.Range("M108").Value = 1167
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H35").Value = 45000000
$ws.Range("I35").Value = 45000000
$ws.Range("K35").Value = 45000000
$ws.Range("M35").Value = -44999702
$ws.Range("H126").Value = 750
$ws.Range("I126").Value = 750
$ws.Range("K126").Value = 2250
$ws.Range("M126").Value = 220
$ws.Range("H132").Value = 1347.7142
$ws.Range("I132").Value = 1239
$ws.Range("K132").Value = 3717
$ws.Range("M132").Value = -1187
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("N7").ClearContents()
$ws.Range("H46").Value = 601
$ws.Range("I46").Value = 601
$ws.Range("K46").Value = 601
$ws.Range("M46").Value = -413
$ws.Range("H122").Value = 12754.2
$ws.Range("I122").Value = 12754.2
$ws.Range("K122").Value = 38262.60000000001
$ws.Range("M122").Value = -35812.60000000001
$ws.Range("H126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("N126").ClearContents()
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 999
$ws.Range("I126").Value = 999
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 2997
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -527
$ws.Range("N126").ClearContents()
$ws.Range("H132").Value = 2637.75
$ws.Range("I132").Value = 2184
$ws.Range("K132").Value = 6552
$ws.Range("M132").Value = -4022
$ws.Range("H136").Value = 3020.6
$ws.Range("I136").Value = 3020.6
$ws.Range("K136").Value = 9061.799999999999
$ws.Range("M136").Value = -6511.799999999999
